$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original (pre-edit) values for the columns that change,
# for rows 4 through 9, before overwriting anything.
$cols = @("A","B","D","E","F","G","H","Q","R","W","Z","AB")

$original = @{}
foreach ($row in 4..9) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$row").Value2
    }
    $original[$row] = $rowData
}

# Mapping of destination row -> source row whose original content it
# should receive (cyclic shift: 4<-5<-7<-9<-6<-8<-4)
$map = @{
    4 = 5
    5 = 7
    7 = 9
    9 = 6
    6 = 8
    8 = 4
}

foreach ($dst in $map.Keys) {
    $src = $map[$dst]
    foreach ($col in $cols) {
        $ws.Range("$col$dst").Value2 = $original[$src][$col]
    }
}
